# #5: property boat&car done
#
# The "汽車" (car) sheet only had an ad-hoc 2-row x 7-col layout (no header
# labels). Bring it in line with the other property sheets: a proper
# header row of field names (adding "capacity" - engine displacement -
# since cars don't have area/share_portion) and the standard trailing
# property_category/category/date/legislator_name/legislator_id/
# source_file/index columns on the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Header row (row 1) ---------------------------------------------------
# Columns B:G already carried header-style (bold/border/center) cells, so
# plain .Value assignment keeps that existing style. The new H:N columns
# need the same look applied explicitly so they pick up the same style.
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"

$newHeaders = @("property_category", "category", "date", "legislator_name", "legislator_id", "source_file", "index")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = 8 + $i
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $newHeaders[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- Data row (row 2) ------------------------------------------------------
$ws.Cells.Item(2, 1).Value = 45
$ws.Cells.Item(2, 2).Value = "BENZS350"
$ws.Cells.Item(2, 3).Value = 3498
$ws.Cells.Item(2, 4).Value = "楊瓊瓔"
$ws.Cells.Item(2, 5).Value = "100年02月25曰"
$ws.Cells.Item(2, 6).Value = "買賣"
$ws.Cells.Item(2, 7).Value = 4720100
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"
# Leading apostrophe forces text so Excel doesn't silently reinterpret the
# ISO-looking "2013-12-25" source-file date stamp as a real date serial.
$ws.Cells.Item(2, 10).Value = "'2013-12-25"
$ws.Cells.Item(2, 11).Value = "楊瓊瓔"
$ws.Cells.Item(2, 12).Value = 854
$ws.Cells.Item(2, 13).Value = "tmp68d11"
$ws.Cells.Item(2, 14).Value = 45
